$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @(2, 68, "house/house006.jpg", "wiegen", "house"),
  @(3, 75, "face/face023.jpg", "sondern", "face"),
  @(4, 53, "face/face022.jpg", "wenden", "face"),
  @(5, 59, "face/face014.jpg", "fühlen", "face"),
  @(6, 3, "face/face016.jpg", "dauern", "face"),
  @(7, 36, "face/face020.jpg", "haken", "face"),
  @(8, 98, "face/face006.jpg", "gründen", "face"),
  @(9, 42, "face/face010.jpg", "scheitern", "face"),
  @(10, 21, "face/face021.jpg", "nehmen", "face"),
  @(11, 74, "face/face017.jpg", "saufen", "face"),
  @(12, 15, "house/house025.jpg", "opfern", "house"),
  @(13, 106, "face/face004.jpg", "mieten", "face"),
  @(14, 110, "house/house001.jpg", "biegen", "house"),
  @(15, 0, "house/house009.jpg", "strahlen", "house"),
  @(16, 121, "house/house023.jpg", "tagen", "house"),
  @(17, 18, "house/house000.jpg", "pflegen", "house"),
  @(18, 100, "face/face027.jpg", "kaufen", "face"),
  @(19, 40, "house/house021.jpg", "loben", "house"),
  @(20, 46, "house/house018.jpg", "krachen", "house"),
  @(21, 81, "house/house019.jpg", "gelten", "house"),
  @(22, 9, "face/face031.jpg", "fesseln", "face"),
  @(23, 80, "house/house031.jpg", "hupen", "house"),
  @(24, 92, "face/face003.jpg", "rasen", "face"),
  @(25, 91, "house/house012.jpg", "laufen", "house"),
  @(26, 124, "house/house004.jpg", "tauschen", "house"),
  @(27, 58, "house/house016.jpg", "schenken", "house"),
  @(28, 26, "face/face001.jpg", "füttern", "face"),
  @(29, 7, "house/house013.jpg", "schätzen", "house"),
  @(30, 94, "house/house010.jpg", "fliegen", "house"),
  @(31, 111, "face/face007.jpg", "hauen", "face"),
  @(32, 2, "house/house028.jpg", "bitten", "house"),
  @(33, 48, "face/face029.jpg", "liefern", "face")
)

foreach ($r in $rows) {
  $ws.Cells.Item($r[0], 2).Value = $r[1]
  $ws.Cells.Item($r[0], 3).Value = $r[2]
  $ws.Cells.Item($r[0], 4).Value = $r[3]
  $ws.Cells.Item($r[0], 5).Value = $r[4]
}